# Refresh the cryptocurrency price/volume snapshot on Sheet1.
# Values that look numeric (e.g. "0.999", "215.41") must be written back as
# literal text -- the source feed renders prices as fixed strings (note the
# "26.991.17"-style thousand-grouped values that aren't valid numbers at
# all), and Excel's normal Value setter would otherwise silently coerce a
# cell like "0.524" into the binary double 0.52400000000000002. Forcing the
# NumberFormat to "@" (Text) before the write keeps the exact digits, and
# ClearFormats() afterwards drops that temporary formatting again so the
# cell's style stays untouched (same default style as every other data
# cell in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.023.81'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +2.22%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.659.58'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  +2.89%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E4'; Value = '  +0.01%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '215.41'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +1.40%  '; ForceText = $false },
    @{ Cell = 'E6'; Value = '  +1.93%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -0.04%  '; ForceText = $false },
    @{ Cell = 'E8'; Value = '  +2.57%  '; ForceText = $false },
    @{ Cell = 'E9'; Value = '  +1.90%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  +4.61%  '; ForceText = $false },
    @{ Cell = 'E11'; Value = '  +4.24%  '; ForceText = $false },
    @{ Cell = 'E12'; Value = '  +2.93%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '1.654.97'; ForceText = $false },
    @{ Cell = 'E13'; Value = '  +2.61%  '; ForceText = $false },
    @{ Cell = 'E14'; Value = '  +2.04%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '0.524'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  +2.79%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '65.58'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  +3.05%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '27.015.03'; ForceText = $false },
    @{ Cell = 'E17'; Value = '  +2.22%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '236.77'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  +1.14%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '0.0₃0739'; ForceText = $false },
    @{ Cell = 'E19'; Value = '  +1.83%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '7.78'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +1.77%  '; ForceText = $false },
    @{ Cell = 'E21'; Value = '  -0.08%  '; ForceText = $false },
    @{ Cell = 'E22'; Value = '  +3.88%  '; ForceText = $false },
    @{ Cell = 'B23'; Value = 'Avalanche'; ForceText = $false },
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; ForceText = $false },
    @{ Cell = 'D23'; Value = '9.30'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  +2.64%  '; ForceText = $false },
    @{ Cell = 'B24'; Value = 'Toncoin'; ForceText = $false },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $false },
    @{ Cell = 'D24'; Value = '2.23'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  +2.06%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '145.19'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -1.07%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '7.14'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  +2.04%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  +0.61%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '15.87'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  +2.46%  '; ForceText = $false },
    @{ Cell = 'E29'; Value = '  -0.03%  '; ForceText = $false },
    @{ Cell = 'E30'; Value = '  +0.20%  '; ForceText = $false },
    @{ Cell = 'E31'; Value = '  +1.49%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '1.560.24'; ForceText = $false },
    @{ Cell = 'E32'; Value = '  +3.45%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '3.31'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  +2.09%  '; ForceText = $false },
    @{ Cell = 'E34'; Value = '  +4.58%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  +8.01%  '; ForceText = $false },
    @{ Cell = 'E36'; Value = '  -0.23%  '; ForceText = $false },
    @{ Cell = 'E37'; Value = '  +3.33%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.901'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +8.89%  '; ForceText = $false },
    @{ Cell = 'E39'; Value = '  +2.79%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '6.05'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  +3.52%  '; ForceText = $false },
    @{ Cell = 'E41'; Value = '  -0.04%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '66.52'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  +8.33%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '0.973'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  +6.42%  '; ForceText = $false },
    @{ Cell = 'E44'; Value = '  +2.69%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '1.801.65'; ForceText = $false },
    @{ Cell = 'E45'; Value = '  +2.94%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '0.775'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +1.67%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '90.26'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  +0.69%  '; ForceText = $false },
    @{ Cell = 'E48'; Value = '  +2.71%  '; ForceText = $false },
    @{ Cell = 'B49'; Value = 'BabyDogeCoin'; ForceText = $false },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; ForceText = $false },
    @{ Cell = 'D49'; Value = '0.0₆0104'; ForceText = $false },
    @{ Cell = 'E49'; Value = '  +1.07%  '; ForceText = $false },
    @{ Cell = 'B50'; Value = 'Algorand'; ForceText = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText = $false },
    @{ Cell = 'D50'; Value = '0.100'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  +4.30%  '; ForceText = $false },
    @{ Cell = 'B51'; Value = 'Cronos'; ForceText = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.0506'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  +1.03%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
